# Auto-generated Excel COM-interop script
# Adds new survey wave column (28. 9. 2021 / aktualizace 6. 10. 2021)
# to sheet "data" (column AI) and sheet "pocetR" (column AH).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "data" ---
$ws1 = $wb.Worksheets.Item("data")

# Header cell AI1: copy formatting from AH1 (bordered, bold, centered header style)
$ws1.Range("AH1").Copy($ws1.Range("AI1"))
$ws1.Range("AI1").Value2 = "28. 9. 2021"

# Data rows 2-67: new values in column AI
$data1 = @{
  2 = 0.09
  3 = 0.18
  4 = 0.73
  5 = 0.07000000000000001
  6 = 0.17
  7 = 0.76
  8 = 0.02
  9 = 0.07000000000000001
  10 = 0.91
  11 = 0.09
  12 = 0.23
  13 = 0.68
  14 = 0.21
  15 = 0.16
  16 = 0.63
  17 = 0.06
  18 = 0.17
  19 = 0.77
  20 = 0.02
  21 = 0.06
  22 = 0.92
  23 = 0.2
  24 = 0.18
  25 = 0.62
  26 = 0.11
  27 = 0.22
  28 = 0.67
  29 = 0.14
  30 = 0.27
  31 = 0.59
  32 = 0.06
  33 = 0.14
  34 = 0.8
  35 = 0.01
  36 = 0.08
  37 = 0.91
  38 = 0.16
  39 = 0.21
  40 = 0.63
  41 = 0.06
  42 = 0.17
  43 = 0.77
  44 = 0.66
  45 = 0.15
  46 = 0.19
  47 = 0.09
  48 = 0.6
  49 = 0.31
  50 = 0.02
  51 = 0.08
  52 = 0.9
  53 = 0.07000000000000001
  54 = 0.19
  55 = 0.74
  56 = 0
  57 = 0.1
  58 = 0.9
  59 = 0.05
  60 = 0.21
  61 = 0.74
  62 = 0.04
  63 = 0.1
  64 = 0.86
  65 = 0.07000000000000001
  66 = 0.09
  67 = 0.84
}
foreach ($row in $data1.Keys) {
  $ws1.Range("AI$row").Value2 = $data1[$row]
}

# Update footer caption in row 68 with new "aktualizace" date
$ws1.Range("A68").Value2 = "Život během pandemie, Zasažení domácností, % respondentů celkově a ve skupinách, aktualizace 6. 10. 2021"

# --- Sheet 2: "pocetR" ---
$ws2 = $wb.Worksheets.Item("pocetR")

# Header cell AH1: copy formatting from AG1
$ws2.Range("AG1").Copy($ws2.Range("AH1"))
$ws2.Range("AH1").Value2 = "28. 9. 2021"

# Data rows 2-23: new values in column AH
$data2 = @{
  2 = 1575
  3 = 750
  4 = 125
  5 = 489
  6 = 211
  7 = 717
  8 = 116
  9 = 95
  10 = 647
  11 = 746
  12 = 525
  13 = 304
  14 = 422
  15 = 1153
  16 = 140
  17 = 283
  18 = 1152
  19 = 274
  20 = 97
  21 = 251
  22 = 137
  23 = 83
}
foreach ($row in $data2.Keys) {
  $ws2.Range("AH$row").Value2 = $data2[$row]
}

# Update footer caption in row 24 with new "aktualizace" date
$ws2.Range("A24").Value2 = "Život během pandemie, Zasažení domácností, velikost dotázaného souboru celkově a ve skupinách, aktualizace 6. 10. 2021"
